$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 69
$ws.Cells.Item($row, 1).Value = 4
$ws.Cells.Item($row, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($row, 3).Value = "Los Lagos"
$ws.Cells.Item($row, 4).Value = 44656
$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($row, 5).Value = 10
$ws.Cells.Item($row, 6).Value = 100112031
$ws.Cells.Item($row, 7).Value = "Poroto verde"
$ws.Cells.Item($row, 8).Value = "Magnum"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 60
$ws.Cells.Item($row, 11).Value = 35000
$ws.Cells.Item($row, 12).Value = 35000
$ws.Cells.Item($row, 13).Value = 35000
$ws.Cells.Item($row, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item($row, 15).Value = "Región Metropolitana"
$ws.Cells.Item($row, 16).Value = 1400
$ws.Cells.Item($row, 17).Value = 25
$ws.Cells.Item($row, 18).Value = "Hortaliza"
